$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows (31 and 32) of Mac-Address / user data following the
# same pattern as the preceding rows.
$newRows = @(
    @{ Row = 31; A = 10001; B = 110030 },
    @{ Row = 32; A = 10001; B = 110031 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = "eng"
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = "superadmin"
    $ws.Cells.Item($row, 6).Value = "now()"
    $ws.Cells.Item($row, 7).Value = "now()"
}

# Update the selected cell/view similar to the post-edit state.
$ws.Range("E28").Select()
